$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Preço médio" (B) values
$ws.Range("B2").Value = 29.56257958202775
$ws.Range("B3").Value = 20.77671525367118
$ws.Range("B4").Value = 20.51043402265566
$ws.Range("B5").Value = 19.50217388325368
$ws.Range("B6").Value = 18.89204421286304
$ws.Range("B7").Value = 16.58840455492738
$ws.Range("B8").Value = 13.70238082987031
$ws.Range("B9").Value = 18.76402348687201

# Update region name in row 7: Piauí -> Alagoas
$ws.Range("A7").Value = "Alagoas"

# Update "Ano" (D) values: 2013-2023 -> 2013-2024
$ws.Range("D2").Value = "2013-2024"
$ws.Range("D3").Value = "2013-2024"
$ws.Range("D4").Value = "2013-2024"
$ws.Range("D5").Value = "2013-2024"
$ws.Range("D6").Value = "2013-2024"
$ws.Range("D7").Value = "2013-2024"
$ws.Range("D8").Value = "2013-2024"
$ws.Range("D9").Value = "2013-2024"
